# This script applies the edits described by the diff:
#  - Updates the "Run Date" cell M1 from 45943 to 45944
#  - For 18 pairs of adjacent rows (same item name in column C), swaps the
#    values of columns B, D, E, F and G between the two rows in each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Run Date cell.
$ws.Range("M1").Value = 45944

# Row pairs whose B/D/E/F/G values need to be swapped.
$rowPairs = @(
    @(316, 317),
    @(346, 347),
    @(351, 352),
    @(355, 356),
    @(372, 373),
    @(379, 380),
    @(389, 390),
    @(419, 420),
    @(421, 422),
    @(457, 458),
    @(581, 582),
    @(583, 584),
    @(586, 587),
    @(601, 602),
    @(687, 688),
    @(709, 710),
    @(720, 721),
    @(872, 873)
)

$cols = @("B", "D", "E", "F", "G")

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")

        $v1 = $cell1.Value()
        $v2 = $cell2.Value()

        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
